$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet has two header rows (row 1 and row 2) describing the
# columns, and 9 data rows (rows 3-11). The new layout merges everything
# into a single, fully-labelled header row (row 1), so the 9 data rows
# shift up by one (rows 2-10).

# Remove the old second header row (row 2); this shifts all data rows up
# by one and keeps their existing values/styles intact.
$ws.Rows("2:2").Delete()

# Reset row 1 formatting/content so we can rebuild the header cleanly.
$ws.Rows("1:1").ClearFormats()
$ws.Rows("1:1").ClearContents()

$ws.Range("A1").Value2 = "idx"
$ws.Range("B1").Value2 = "idx2"
$ws.Range("C1").Value2 = "Name"
$ws.Range("D1").Value2 = "Date Start"
$ws.Range("E1").Value2 = "Date End"
$ws.Range("F1").Value2 = "(m3/s)"
$ws.Range("G1").Value2 = "(MW1)"
$ws.Range("H1").Value2 = "(MW2)"
$ws.Range("I1").Value2 = "(GWh) Winter"
$ws.Range("J1").Value2 = "(GWh) Summer"
$ws.Range("K1").Value2 = "(GWh) Year"

# F1:K1 use the smaller 9pt data font (matches the rest of the table),
# while A1:E1 keep the default 10pt font.
$ws.Range("F1:K1").Font.Size = 9

# Match the selection left behind by the edit (row 2, the first data row).
$ws.Range("A2:K2").Select() | Out-Null
